# Updated cryptos list on Wed Apr  3 18:35:31 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.618.24"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.289.09"
$ws.Range("E3").Value = "  +0.65%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "184.71"
$ws.Range("E5").Value = "  +0.71%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "551.73"
$ws.Range("E6").Value = "  -0.57%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.21%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.281.49"
$ws.Range("E8").Value = "  +0.55%  "

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  -3.04%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -6.07%  "

# Row 11 - Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.574"
$ws.Range("E11").Value = "  -1.92%  "

# Row 12 - Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.47"

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -1.65%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.823.07"
$ws.Range("E14").Value = "  +0.64%  "

# Row 15 - Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.39"
$ws.Range("E15").Value = "  -1.86%  "

# Row 16 - BitcoinCash
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "576.60"
$ws.Range("E16").Value = "  -8.40%  "

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.687.63"
$ws.Range("E17").Value = "  +0.10%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  +0.43%  "

# Row 19 - WrappedEther
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.308.16"
$ws.Range("E19").Value = "  +1.05%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.61"
$ws.Range("E20").Value = "  -1.24%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.77"
$ws.Range("E21").Value = "  -4.91%  "

# Row 22 - Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.886"
$ws.Range("E22").Value = "  -1.90%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.90"
$ws.Range("E23").Value = "  +2.06%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +1.01%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.70"
$ws.Range("E25").Value = "  -8.13%  "

# Row 26 - PancakeSwap
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.91"
$ws.Range("E26").Value = "  -1.50%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  +0.34%  "

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.31"
$ws.Range("E28").Value = "  -2.41%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -0.19%  "

# Row 30 - Filecoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.31"
$ws.Range("E30").Value = "  -3.86%  "

# Row 31 - NEARProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.56"
$ws.Range("E31").Value = "  +4.38%  "

# Row 32 - Bittensor
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "566.56"
$ws.Range("E32").Value = "  +4.58%  "

# Row 33 - dogwifhat
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.65"
$ws.Range("E33").Value = "  -9.31%  "

# Row 34 - Cosmos
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.77"
$ws.Range("E34").Value = "  -1.96%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  -2.16%  "

# Row 36 - Maker
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.706.60"
$ws.Range("E36").Value = "  +0.55%  "

# Row 37 - Dai
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.18%  "

# Row 38 - OKB
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.18"
$ws.Range("E38").Value = "  -3.62%  "

# Row 39 - InjectiveProtocol
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.56"
$ws.Range("E39").Value = "  +3.65%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -3.93%  "

# Row 41 - PEPE
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0681"
$ws.Range("E41").Value = "  -5.42%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -8.04%  "

# Row 43 / 44 - ApeXProtocol and Fetch.AI swap ranking positions
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("E43").Value = "  -4.93%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.34"
$ws.Range("E44").Value = "  +3.06%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  -1.53%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  -2.09%  "

# Row 47 - CoreDAO
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  -12.30%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -1.02%  "

# Row 49 - FirstDigitalUSD
$ws.Range("E49").Value = "  +0.07%  "

# Row 50 - ThetaToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.49"
$ws.Range("E50").Value = "  -4.33%  "

# Row 51 - Monero
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "127.46"
$ws.Range("E51").Value = "  +5.54%  "
